# Copy edit pass over "ADS-Geo Script.docx" per commit "Geo Script for
# Usability Testing / copy edit".

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# "...observe your interaction with the site as well." -> drop "as well"
Replace-Text `
    "I have a colleague helping me take notes and observe your interaction with the site as well." `
    "I have a colleague helping me take notes and observe your interaction with the site."

# Oxford comma: "comments or areas" -> "comments, or areas"
Replace-Text `
    "If you have any questions, comments or areas of confusion while you are working, please let me know." `
    "If you have any questions, comments, or areas of confusion while you are working, please let me know."

# Grammar fix: "ask that you to try work through the tasks based on what you see
# on screen, but if you reach" -> "ask that you try to work through the tasks
# based on what you see on the screen. If you reach"
Replace-Text `
    "I would ask that you to try work through the tasks based on what you see on screen, but if you reach a point where you are not sure where or how to find something, please feel free to use ‘Search’." `
    "I would ask that you try to work through the tasks based on what you see on the screen. If you reach a point where you are not sure where or how to find something, please feel free to use ‘Search’."

# "for reference if needed." -> "for reference."
Replace-Text `
    "We will be recording this session for reference if needed. We are capturing your" `
    "We will be recording this session for reference. We are capturing your"

# "voice and what you see" -> "voice, and what you see"
Replace-Text `
    " voice and what you see on the screen. Your name will not be associated or reported with data or findings from this evaluation. " `
    " voice, and what you see on the screen. Your name will not be associated or reported with data or findings from this evaluation. "

# Scenario 1 copy edits
Replace-Text `
    "Your own mother has called you complaining of X symptoms. After a brief discussion, you have learned she has recently began taking a new prescription given to her by her physician. " `
    "Your mother has called you complaining of X symptoms. After a brief discussion, you realize she has started taking a new prescription given to her by her physician. "

Replace-Text `
    "Your mother takes a number of other prescriptions for various symptoms associated with the elderly as well as a variety of OTC medications and supplements. " `
    "Your mother takes a number of other prescriptions for various symptoms associated with the elderly, as well as a variety of over the counter (OTC) medications and supplements. "

Replace-Text `
    "How would you determine if this new drug could have an unintended interaction with some of her other routine medications?" `
    "How would you determine if this new drug could have an unintended interaction with some of the other routine medications?"

# Scenario 2 copy edits: Ibuprofin -> Ibuprofen (both occurrences, same fix)
Replace-Text `
    "Your Allergies have been flaring up. You’ve been taking both Benadryl and Ibuprofin. However, the sinus pressure is still causing some mild pain. " `
    "Your Allergies have been flaring up. You’ve been taking both Benadryl and Ibuprofen. However, the sinus pressure is still causing some mild pain. "

Replace-Text `
    "After determining that the Ibuprofin didn’t work, you consider taking Aspirin instead. " `
    "After determining that the Ibuprofen didn’t work, you consider taking Aspirin instead. "

Replace-Text `
    "if taking these three drugs together in such a short time period would be advised?" `
    "if taking these three drugs together in such a short time period would cause adverse reactions?"

# Scenario 3 copy edit
Replace-Text `
    "Your doctor has prescribed you a medication. You have never been prescribed any OTC medication before. " `
    "Your doctor has prescribed a medication for you to take. You have never been prescribed any OTC medication before. "

# Move the hidden "_GoBack" bookmark from the end of the document (after the
# Scenario 3 intro paragraphs) to just before the "?" that now ends the
# Scenario 2 closing question -- this is where Word leaves it after the
# edits above land as the most-recent-edit location.
$q = $d.Content.Find.Execute("cause adverse reactions?", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
$anchor = $d.Content
$anchor.Find.Execute("cause adverse reactions?", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$bmPos = $anchor.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
